# Update cryptocurrency price/volume data per Oct 9 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.154.60"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.447.22"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.75"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.82"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.442.03"
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.40"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.873.51"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.129.08"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.440.62"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.74"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.10"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -6.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.63"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.07"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "599.17"
$ws.Range("E27").Value = "  -5.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0970"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.568.43"
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.43"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "153.04"
$ws.Range("E39").Value = "  +4.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.41"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.29"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.20"
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.71"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.53"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "141.82"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0268"
$ws.Range("E48").Value = "  +20.24%  "
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.90"
$ws.Range("E51").Value = "  +0.85%  "
